# Refresh LCII yearly financials: add FY2018 column and correct a few
# previously-reported prior-year figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D for the new fiscal year (2018-12-31),
# shifting the existing year columns (old D:K) right to E:L.
$ws.Range("D5:D102").EntireColumn.Insert()

# The newly inserted column inherits formatting from the column to its left (C);
# copy the number/date formatting from column E (the old column D) onto it so
# dates and numbers render the same as the rest of the table.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate new column D with FY2018 figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2475800
$ws.Range("D9").Value = 1955500
$ws.Range("D10").Value = 520300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 26500
$ws.Range("D17").Value = 2277000
$ws.Range("D18").Value = 198800
$ws.Range("D20").Value = -6400
$ws.Range("D21").Value = 259900
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 192400
$ws.Range("D24").Value = 43200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 149200
$ws.Range("D27").Value = 149200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -600
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 6400
$ws.Range("D33").Value = 148600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 148600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 14900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 121800
$ws.Range("D44").Value = 340600
$ws.Range("D45").Value = 49300
$ws.Range("D46").Value = 526700
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 322900
$ws.Range("D49").Value = 356500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 37900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1243900
$ws.Range("D57").Value = 78400
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 99200
$ws.Range("D60").Value = 177600
$ws.Range("D61").Value = 293500
$ws.Range("D62").Value = 66500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 537600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 563500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 706300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 148600
$ws.Range("D83").Value = 67500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 156600
$ws.Range("D91").Value = -119800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -302800
$ws.Range("D96").Value = -59300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 135100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -11100

# A handful of previously-reported figures were corrected as part of this refresh
$ws.Range("E24").Value = 66800
$ws.Range("E26").Value = 146100
$ws.Range("E27").Value = 146100
$ws.Range("E29").Value = -13200
$ws.Range("E89").Value = 152700
$ws.Range("F89").Value = 201700
$ws.Range("E100").Value = -66900
$ws.Range("F100").Value = -36100
